$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 280.2
$ws.Cells.Item(5, 9).Value = 325.25
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 325.25
$ws.Cells.Item(5, 12).Value = 100
$ws.Cells.Item(5, 13).Value = -210.25
$ws.Cells.Item(5, 14).Value = -330

$ws.Cells.Item(32, 8).Value = 1234.2858
$ws.Cells.Item(32, 9).Value = 879.3333
$ws.Cells.Item(32, 10).Value = 1500.5
$ws.Cells.Item(32, 11).Value = 879.3333
$ws.Cells.Item(32, 12).Value = 1500.5
$ws.Cells.Item(32, 13).Value = -553.3333
$ws.Cells.Item(32, 14).Value = -2152.5

$ws.Cells.Item(121, 8).Value = 2342.7144
$ws.Cells.Item(121, 10).Value = 2566.5
$ws.Cells.Item(121, 12).Value = 7699.5
$ws.Cells.Item(121, 14).Value = -11193.5

$ws.Cells.Item(128, 8).Value = 22983.334
$ws.Cells.Item(128, 10).Value = 22983.334
$ws.Cells.Item(128, 12).Value = 22983.334
$ws.Cells.Item(128, 14).Value = -32943.334

$ws.Cells.Item(132, 8).Value = 1721.44
$ws.Cells.Item(132, 9).Value = 1214.2632
$ws.Cells.Item(132, 10).Value = 3327.5
$ws.Cells.Item(132, 11).Value = 3642.7896
$ws.Cells.Item(132, 12).Value = 9982.5
$ws.Cells.Item(132, 13).Value = -1112.7896
$ws.Cells.Item(132, 14).Value = -15042.5

$ws.Cells.Item(137, 8).Value = 1308
$ws.Cells.Item(137, 9).Value = 1120.0714
$ws.Cells.Item(137, 10).Value = 1600.3334
$ws.Cells.Item(137, 11).Value = 3360.2142
$ws.Cells.Item(137, 12).Value = 4801.0002
$ws.Cells.Item(137, 13).Value = -810.2142000000003
$ws.Cells.Item(137, 14).Value = -9901.0002

$ws.Cells.Item(138, 8).Value = 2901625.2
$ws.Cells.Item(138, 10).Value = 3799.7874
$ws.Cells.Item(138, 12).Value = 11399.3622
$ws.Cells.Item(138, 14).Value = -21679.3622

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(57, 8).Value = 3500
$ws.Cells.Item(57, 9).Value = 3500
$ws.Cells.Item(57, 11).Value = 3500
$ws.Cells.Item(57, 13).Value = -3016

$ws.Cells.Item(74, 8).Value = 1136
$ws.Cells.Item(74, 9).Value = 964.6842
$ws.Cells.Item(74, 10).Value = 1386.3846
$ws.Cells.Item(74, 11).Value = 964.6842
$ws.Cells.Item(74, 12).Value = 1386.3846
$ws.Cells.Item(74, 13).Value = -90.68420000000003
$ws.Cells.Item(74, 14).Value = -3134.3846

$ws.Cells.Item(77, 8).Value = 1136
$ws.Cells.Item(77, 9).Value = 964.6842
$ws.Cells.Item(77, 10).Value = 1386.3846
$ws.Cells.Item(77, 11).Value = 4823.421
$ws.Cells.Item(77, 12).Value = 6931.923000000001
$ws.Cells.Item(77, 13).Value = -455.4210000000003
$ws.Cells.Item(77, 14).Value = -15667.923

$ws.Cells.Item(102, 8).Value = 1576.6666
$ws.Cells.Item(102, 9).Value = 1615
$ws.Cells.Item(102, 11).Value = 1615
$ws.Cells.Item(102, 13).Value = 7

$ws.Cells.Item(132, 8).Value = 2084.1052
$ws.Cells.Item(132, 9).Value = 1457.1428
$ws.Cells.Item(132, 10).Value = 3839.6
$ws.Cells.Item(132, 11).Value = 4371.428400000001
$ws.Cells.Item(132, 12).Value = 11518.8
$ws.Cells.Item(132, 13).Value = -1841.428400000001
$ws.Cells.Item(132, 14).Value = -16578.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1751.5
$ws.Cells.Item(94, 9).Value = 1303
$ws.Cells.Item(94, 10).Value = 2200
$ws.Cells.Item(94, 11).Value = 1303
$ws.Cells.Item(94, 12).Value = 2200
$ws.Cells.Item(94, 13).Value = -852
$ws.Cells.Item(94, 14).Value = -3102

$ws.Cells.Item(107, 8).Value = 52701.2
$ws.Cells.Item(107, 9).Value = 65489
$ws.Cells.Item(107, 10).Value = 1550
$ws.Cells.Item(107, 11).Value = 65489
$ws.Cells.Item(107, 12).Value = 1550
$ws.Cells.Item(107, 13).Value = -63569
$ws.Cells.Item(107, 14).Value = -5390

$ws.Cells.Item(122, 8).Value = 63598.184
$ws.Cells.Item(122, 10).Value = 63598.184
$ws.Cells.Item(122, 12).Value = 63598.184
$ws.Cells.Item(122, 14).Value = -73398.18400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1725.697
$ws.Cells.Item(31, 9).Value = 1251
$ws.Cells.Item(31, 10).Value = 2817.5
$ws.Cells.Item(31, 11).Value = 1251
$ws.Cells.Item(31, 12).Value = 2817.5
$ws.Cells.Item(31, 13).Value = -956
$ws.Cells.Item(31, 14).Value = -3407.5

$ws.Cells.Item(34, 8).Value = 1725.697
$ws.Cells.Item(34, 9).Value = 1251
$ws.Cells.Item(34, 10).Value = 2817.5
$ws.Cells.Item(34, 11).Value = 1251
$ws.Cells.Item(34, 12).Value = 2817.5
$ws.Cells.Item(34, 13).Value = -1049
$ws.Cells.Item(34, 14).Value = -3221.5

$ws.Cells.Item(94, 8).Value = 1844.238
$ws.Cells.Item(94, 9).Value = 1050
$ws.Cells.Item(94, 10).Value = 2031.1177
$ws.Cells.Item(94, 11).Value = 1050
$ws.Cells.Item(94, 12).Value = 2031.1177
$ws.Cells.Item(94, 13).Value = -599
$ws.Cells.Item(94, 14).Value = -2933.1177

$ws.Cells.Item(109, 8).Value = 30285
$ws.Cells.Item(109, 10).Value = 30285
$ws.Cells.Item(109, 12).Value = 30285
$ws.Cells.Item(109, 14).Value = -32365

$ws.Cells.Item(112, 8).Value = 30327
$ws.Cells.Item(112, 10).Value = 30327
$ws.Cells.Item(112, 12).Value = 30327
$ws.Cells.Item(112, 14).Value = -33281

$ws.Cells.Item(132, 8).Value = 2355.3635
$ws.Cells.Item(132, 9).Value = 1767.8889
$ws.Cells.Item(132, 11).Value = 5303.6667
$ws.Cells.Item(132, 13).Value = -2773.6667

$ws.Cells.Item(141, 8).Value = 32722.25
$ws.Cells.Item(141, 10).Value = 30296.334
$ws.Cells.Item(141, 12).Value = 30296.334
$ws.Cells.Item(141, 14).Value = -40656.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 741.1111
$ws.Cells.Item(107, 9).Value = 1004
$ws.Cells.Item(107, 10).Value = 412.5
$ws.Cells.Item(107, 11).Value = 3012
$ws.Cells.Item(107, 12).Value = 1237.5
$ws.Cells.Item(107, 13).Value = -1092
$ws.Cells.Item(107, 14).Value = -5077.5

$ws.Cells.Item(131, 8).Value = 21278268
$ws.Cells.Item(131, 9).Value = 515
$ws.Cells.Item(131, 10).Value = 22223946
$ws.Cells.Item(131, 11).Value = 1545
$ws.Cells.Item(131, 12).Value = 66671838
$ws.Cells.Item(131, 13).Value = 3495
$ws.Cells.Item(131, 14).Value = -66681918

$ws.Cells.Item(136, 8).Value = 3515.0417
$ws.Cells.Item(136, 10).Value = 9476.143
$ws.Cells.Item(136, 12).Value = 28428.429
$ws.Cells.Item(136, 14).Value = -38628.429

$ws.Cells.Item(137, 8).Value = 3190.4285
$ws.Cells.Item(137, 9).Value = 2000
$ws.Cells.Item(137, 10).Value = 3388.8333
$ws.Cells.Item(137, 11).Value = 6000
$ws.Cells.Item(137, 12).Value = 10166.4999
$ws.Cells.Item(137, 13).Value = -900
$ws.Cells.Item(137, 14).Value = -20366.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4095.7
$ws.Cells.Item(132, 9).Value = 3315.3333
$ws.Cells.Item(132, 10).Value = 4430.143
$ws.Cells.Item(132, 11).Value = 9945.999899999999
$ws.Cells.Item(132, 12).Value = 13290.429
$ws.Cells.Item(132, 13).Value = -7415.999899999999
$ws.Cells.Item(132, 14).Value = -18350.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4150
$ws.Cells.Item(7, 9).Value = 5266.6665
$ws.Cells.Item(7, 10).Value = 3033.3333
$ws.Cells.Item(7, 11).Value = 5266.6665
$ws.Cells.Item(7, 12).Value = 3033.3333
$ws.Cells.Item(7, 13).Value = -5154.6665
$ws.Cells.Item(7, 14).Value = -3257.3333

$ws.Cells.Item(61, 8).Value = 18241.084
$ws.Cells.Item(61, 9).Value = 25923.5
$ws.Cells.Item(61, 10).Value = 2876.25
$ws.Cells.Item(61, 11).Value = 25923.5
$ws.Cells.Item(61, 12).Value = 2876.25
$ws.Cells.Item(61, 13).Value = -25721.5
$ws.Cells.Item(61, 14).Value = -3280.25

$ws.Cells.Item(98, 8).Value = 21052
$ws.Cells.Item(98, 10).Value = 21052
$ws.Cells.Item(98, 12).Value = 21052
$ws.Cells.Item(98, 14).Value = -27042

$ws.Cells.Item(113, 8).Value = 18241.084
$ws.Cells.Item(113, 9).Value = 25923.5
$ws.Cells.Item(113, 10).Value = 2876.25
$ws.Cells.Item(113, 11).Value = 25923.5
$ws.Cells.Item(113, 12).Value = 2876.25
$ws.Cells.Item(113, 13).Value = -23753.5
$ws.Cells.Item(113, 14).Value = -7216.25

$ws.Cells.Item(122, 8).Value = 19237884
$ws.Cells.Item(122, 9).Value = 50005800
$ws.Cells.Item(122, 10).Value = 7936.5
$ws.Cells.Item(122, 11).Value = 150017400
$ws.Cells.Item(122, 12).Value = 23809.5
$ws.Cells.Item(122, 13).Value = -150014950
$ws.Cells.Item(122, 14).Value = -28709.5

$ws.Cells.Item(126, 8).Value = 4150
$ws.Cells.Item(126, 9).Value = 5266.6665
$ws.Cells.Item(126, 10).Value = 3033.3333
$ws.Cells.Item(126, 11).Value = 15799.9995
$ws.Cells.Item(126, 12).Value = 9099.999899999999
$ws.Cells.Item(126, 13).Value = -13329.9995
$ws.Cells.Item(126, 14).Value = -14039.9999

$ws.Cells.Item(132, 8).Value = 6035.8696
$ws.Cells.Item(132, 9).Value = 5836.769
$ws.Cells.Item(132, 11).Value = 17510.307
$ws.Cells.Item(132, 13).Value = -14980.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 19234086
$ws.Cells.Item(122, 9).Value = 31252888
$ws.Cells.Item(122, 11).Value = 93758664
$ws.Cells.Item(122, 13).Value = -93756214

$ws.Cells.Item(132, 8).Value = 2231.8125
$ws.Cells.Item(132, 9).Value = 1668.4445
$ws.Cells.Item(132, 11).Value = 5005.333500000001
$ws.Cells.Item(132, 13).Value = -2475.333500000001
